$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$data = @(
    @(2, 5, 6),
    @(3, 4, 5),
    @(4, 5, 6),
    @(5, 7, 7),
    @(6, 3, 4),
    @(7, 7, 7),
    @(8, 4, 5),
    @(9, 6, 6),
    @(10, 8, 8),
    @(11, 8, 8),
    @(12, 11, 11),
    @(13, 7, 7),
    @(14, 5, 7),
    @(15, 6, 7),
    @(16, 5, 6),
    @(17, 8, 8),
    @(18, 7, 7),
    @(19, 7, 7),
    @(20, 7, 7),
    @(21, 6, 6),
    @(22, 7, 7),
    @(23, 6, 6),
    @(24, 6, 7),
    @(25, 5, 5),
    @(26, 7, 7),
    @(27, 9, 9),
    @(28, 4, 5),
    @(29, 9, 9),
    @(30, 6, 6),
    @(31, 9, 9),
    @(32, 6, 7),
    @(33, 6, 7),
    @(34, 5, 5),
    @(35, 5, 5),
    @(36, 6, 6),
    @(37, 8, 9),
    @(38, 6, 7),
    @(39, 7, 7),
    @(40, 6, 6),
    @(41, 5, 6),
    @(42, 9, 9),
    @(43, 5, 5),
    @(44, 5, 6),
    @(45, 7, 8),
    @(46, 9, 9),
    @(47, 9, 9),
    @(48, 9, 9),
    @(49, 9, 9),
    @(50, 6, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}
